$p = $ppt.ActivePresentation

# Slide 1: Title "First" + " " + "slide" -> single run "First slide"
# (Re-assigning the identical text is a no-op for the engine's run-consolidation,
# so force a transient change first.)
$s1 = $p.Slides.Item(1)
$t1 = $s1.Shapes.Item(1).TextFrame.TextRange
$t1.Text = "__tmp__"
$t1.Text = "First slide"

# Slide 3: Title "Third" + " " + "slide" -> single run "Third slide"
$s3 = $p.Slides.Item(3)
$t3 = $s3.Shapes.Item(1).TextFrame.TextRange
$t3.Text = "__tmp__"
$t3.Text = "Third slide"

# Slide 2's notes page: consolidate the many single-word runs into one run.
$s2 = $p.Slides.Item(2)
$notesPage = $s2.NotesPage
$tn = $notesPage.Shapes.Item(2).TextFrame.TextRange
$tn.Text = "__tmp__"
$tn.Text = "Some notes here: this first slide should use the Blank template"
